$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 24,14

$arr[0,0] = 0.5576086467745256
$arr[0,1] = 0.1401259805752488
$arr[0,2] = 0
$arr[0,3] = 0.2158203427384855
$arr[0,4] = 2.101583851666291
$arr[0,5] = 0.002480846746097597
$arr[0,6] = 0
$arr[0,7] = 1.030516627725682
$arr[0,8] = 0.06832021374958686
$arr[0,9] = 0.2810406009705275
$arr[0,10] = 0.4458249572567325
$arr[0,11] = 0.2357755197726483
$arr[0,12] = 2.090795578000424
$arr[0,13] = 3.870837430261133
$arr[1,0] = 0.524764461952401
$arr[1,1] = 0.1394233556882121
$arr[1,2] = 0
$arr[1,3] = 0.2156799295414693
$arr[1,4] = 2.103628181986785
$arr[1,5] = 0.002482973314430057
$arr[1,6] = 0
$arr[1,7] = 1.037747933989987
$arr[1,8] = 0.06647462156653106
$arr[1,9] = 0.2522360085775546
$arr[1,10] = 0.4409351929450196
$arr[1,11] = 0.2282478055067294
$arr[1,12] = 2.109299138222191
$arr[1,13] = 3.894317930062471
$arr[2,0] = 0.5047659655288896
$arr[2,1] = 0.1389879205728768
$arr[2,2] = 0
$arr[2,3] = 0.215661839952233
$arr[2,4] = 2.105797412865122
$arr[2,5] = 0.002484349641552898
$arr[2,6] = 0
$arr[2,7] = 1.042590932208917
$arr[2,8] = 0.06532875163213703
$arr[2,9] = 0.2345743598783798
$arr[2,10] = 0.4381004256602097
$arr[2,11] = 0.2237133670508449
$arr[2,12] = 2.121242743642517
$arr[2,13] = 3.91039846181512
$arr[3,0] = 0.496659307014653
$arr[3,1] = 0.1388094739337262
$arr[3,2] = 0
$arr[3,3] = 0.2156716530205323
$arr[3,4] = 2.106911599114866
$arr[3,5] = 0.00248492831153518
$arr[3,6] = 0
$arr[3,7] = 1.044665907256444
$arr[3,8] = 0.06485863060783004
$arr[3,9] = 0.2273837102862615
$arr[3,10] = 0.4369875523157702
$arr[3,11] = 0.221887753753137
$arr[3,12] = 2.126256364782337
$arr[3,13] = 3.917370019741384
$arr[4,0] = 0.4953158123266519
$arr[4,1] = 0.1387797826359574
$arr[4,2] = 0
$arr[4,3] = 0.2156743219814317
$arr[4,4] = 2.107110523807918
$arr[4,5] = 0.002485025476326727
$arr[4,6] = 0
$arr[4,7] = 1.045016583076638
$arr[4,8] = 0.06478037639372403
$arr[4,9] = 0.2261901214485675
$arr[4,10] = 0.4368053217295511
$arr[4,11] = 0.2215859585293423
$arr[4,12] = 2.127097723375996
$arr[4,13] = 3.918552934485419
$arr[5,0] = 0.5046564617648244
$arr[5,1] = 0.1389855180318023
$arr[5,2] = 0
$arr[5,3] = 0.2156619026422355
$arr[5,4] = 2.105811506575058
$arr[5,5] = 0.00248435737354657
$arr[5,6] = 0
$arr[5,7] = 1.042618505307754
$arr[5,8] = 0.06532242422886014
$arr[5,9] = 0.234477356869732
$arr[5,10] = 0.4380852454793569
$arr[5,11] = 0.2236886560402311
$arr[5,12] = 2.121309765804188
$arr[5,13] = 3.910490787303502
$arr[6,0] = 0.5462494924468331
$arr[6,1] = 0.1398845544080238
$arr[6,2] = 0
$arr[6,3] = 0.2157578211045568
$arr[6,4] = 2.102099255172988
$arr[6,5] = 0.002481565363662961
$arr[6,6] = 0
$arr[6,7] = 1.032926404222721
$arr[6,8] = 0.06768649212877165
$arr[6,9] = 0.2711039610790351
$arr[6,10] = 0.4441042941278255
$arr[6,11] = 0.2331618843857406
$arr[6,12] = 2.097054783464357
$arr[6,13] = 3.878588463447329
$arr[7,0] = 0.6291215714020382
$arr[7,1] = 0.1416153998231309
$arr[7,2] = 0
$arr[7,3] = 0.2164842427352234
$arr[7,4] = 2.10205720785693
$arr[7,5] = 0.002476648124384232
$arr[7,6] = 0
$arr[7,7] = 1.017113571594287
$arr[7,8] = 0.07222138362499209
$arr[7,9] = 0.3431062973053258
$arr[7,10] = 0.4572302069445158
$arr[7,11] = 0.2524272946668447
$arr[7,12] = 2.054107661792014
$arr[7,13] = 3.829213029403292
$arr[8,0] = 0.6907787100357439
$arr[8,1] = 0.1428672103891699
$arr[8,2] = 0
$arr[8,3] = 0.2173433965643206
$arr[8,4] = 2.106421786848244
$arr[8,5] = 0.002473372252285889
$arr[8,6] = 0
$arr[8,7] = 1.007437576894006
$arr[8,8] = 0.07549115100110271
$arr[8,9] = 0.3960975559208464
$arr[8,10] = 0.4676720957807561
$arr[8,11] = 0.266993893269003
$arr[8,12] = 2.025362697930817
$arr[8,13] = 3.800958330527322
$arr[9,0] = 0.7189903017925587
$arr[9,1] = 0.1434323409223381
$arr[9,2] = 0
$arr[9,3] = 0.217804354553369
$arr[9,4] = 2.109358397503613
$arr[9,5] = 0.002471954414296605
$arr[9,6] = 0
$arr[9,7] = 1.003456289743909
$arr[9,8] = 0.07696510957072888
$arr[9,9] = 0.4202210982752206
$arr[9,10] = 0.4725941069750519
$arr[9,11] = 0.273708670113848
$arr[9,12] = 2.012894329930683
$arr[9,13] = 3.789842994889682
$arr[10,0] = 0.7296962156747497
$arr[10,1] = 0.1436457131749904
$arr[10,2] = 0
$arr[10,3] = 0.2179889397021881
$arr[10,4] = 2.110606856182173
$arr[10,5] = 0.002471427870641427
$arr[10,6] = 0
$arr[10,7] = 1.0020090469853
$arr[10,8] = 0.07752130797511825
$arr[10,9] = 0.4293581789051473
$arr[10,10] = 0.4744825042968301
$arr[10,11] = 0.2762639273348313
$arr[10,12] = 2.008260227103541
$arr[10,13] = 3.785883526621092
$arr[11,0] = 0.7273895029633763
$arr[11,1] = 0.1435997878167754
$arr[11,2] = 0
$arr[11,3] = 0.2179487407360661
$arr[11,4] = 2.110331916101785
$arr[11,5] = 0.002471540811248304
$arr[11,6] = 0
$arr[11,7] = 1.002318052026407
$arr[11,8] = 0.07740160808793206
$arr[11,9] = 0.4273902633782711
$arr[11,10] = 0.4740747157662923
$arr[11,11] = 0.2757130532865304
$arr[11,12] = 2.009254377196978
$arr[11,13] = 3.786725168640913
$arr[12,0] = 0.7198706309526983
$arr[12,1] = 0.1434499078697868
$arr[12,2] = 0
$arr[12,3] = 0.217819339767491
$arr[12,4] = 2.10945837671035
$arr[12,5] = 0.002471910887893777
$arr[12,6] = 0
$arr[12,7] = 1.003336014389124
$arr[12,8] = 0.07701090769903374
$arr[12,9] = 0.4209727738363824
$arr[12,10] = 0.4727489758419807
$arr[12,11] = 0.2739186428870397
$arr[12,12] = 2.012511327985114
$arr[12,13] = 3.7895122446767
$arr[13,0] = 0.7152680524962705
$arr[13,1] = 0.143358019722271
$arr[13,2] = 0
$arr[13,3] = 0.2177413826004795
$arr[13,4] = 2.108941065302432
$arr[13,5] = 0.002472138918436828
$arr[13,6] = 0
$arr[13,7] = 1.003967407633152
$arr[13,8] = 0.0767713368792613
$arr[13,9] = 0.4170421240765165
$arr[13,10] = 0.4719401119484274
$arr[13,11] = 0.2728211391750435
$arr[13,12] = 2.014517685812045
$arr[13,13] = 3.791251916409379
$arr[14,0] = 0.6889382447231753
$arr[14,1] = 0.1428301900272402
$arr[14,2] = 0
$arr[14,3] = 0.2173146783908457
$arr[14,4] = 2.106248982360114
$arr[14,5] = 0.002473466363575595
$arr[14,6] = 0
$arr[14,7] = 1.007706216788304
$arr[14,8] = 0.07539455148375396
$arr[14,9] = 0.3945213357115165
$arr[14,10] = 0.4673538747579187
$arr[14,11] = 0.2665568286640294
$arr[14,12] = 2.02618976571638
$arr[14,13] = 3.801719694082891
$arr[15,0] = 0.6728271206287957
$arr[15,1] = 0.1425052699118439
$arr[15,2] = 0
$arr[15,3] = 0.2170708317225056
$arr[15,4] = 2.104840841619904
$arr[15,5] = 0.002474299211515207
$arr[15,6] = 0
$arr[15,7] = 1.010107476219396
$arr[15,8] = 0.07454647277822346
$arr[15,9] = 0.3807096977678555
$arr[15,10] = 0.4645842728560439
$arr[15,11] = 0.2627363733977859
$arr[15,12] = 2.033505886874179
$arr[15,13] = 3.808586278846775
$arr[16,0] = 0.6635758551883839
$arr[16,1] = 0.1423179779334376
$arr[16,2] = 0
$arr[16,3] = 0.2169371812242709
$arr[16,4] = 2.104120462954995
$arr[16,5] = 0.002474785058267035
$arr[16,6] = 0
$arr[16,7] = 1.011528188316152
$arr[16,8] = 0.07405741321568371
$arr[16,9] = 0.3727672913364302
$arr[16,10] = 0.4630074664230932
$arr[16,11] = 0.2605472765564798
$arr[16,12] = 2.037771141489469
$arr[16,13] = 3.81269934721962
$arr[17,0] = 0.6604462068303292
$arr[17,1] = 0.1422544945297872
$arr[17,2] = 0
$arr[17,3] = 0.216893065275606
$arr[17,4] = 2.103891945341985
$arr[17,5] = 0.002474950729688373
$arr[17,6] = 0
$arr[17,7] = 1.012016015384919
$arr[17,8] = 0.07389160902829417
$arr[17,9] = 0.3700784333379943
$arr[17,10] = 0.462476374108121
$arr[17,11] = 0.259807522068968
$arr[17,12] = 2.039225110597546
$arr[17,13] = 3.81412006549823
$arr[18,0] = 0.6745405857776348
$arr[18,1] = 0.1425399003478987
$arr[18,2] = 0
$arr[18,3] = 0.2170961065506098
$arr[18,4] = 2.104981475361413
$arr[18,5] = 0.002474209848459994
$arr[18,6] = 0
$arr[18,7] = 1.009847763053457
$arr[18,8] = 0.07463688349090347
$arr[18,9] = 0.3821797997415786
$arr[18,10] = 0.4648774269038114
$arr[18,11] = 0.263142206835596
$arr[18,12] = 2.03272115182001
$arr[18,13] = 3.807838389885262
$arr[19,0] = 0.722078492875994
$arr[19,1] = 0.1434939484211739
$arr[19,2] = 0
$arr[19,3] = 0.217857076173658
$arr[19,4] = 2.109711256543164
$arr[19,5] = 0.002471801906518887
$arr[19,6] = 0
$arr[19,7] = 1.003035375868784
$arr[19,8] = 0.0771257191662329
$arr[19,9] = 0.4228576953205163
$arr[19,10] = 0.4731377133161487
$arr[19,11] = 0.274445366517611
$arr[19,12] = 2.011552310727732
$arr[19,13] = 3.788686838836725
$arr[20,0] = 0.753279809234698
$arr[20,1] = 0.1441137954980505
$arr[20,2] = 0
$arr[20,3] = 0.218412853636206
$arr[20,4] = 2.113597487795673
$arr[20,5] = 0.002470288546484109
$arr[20,6] = 0
$arr[20,7] = 0.9989350230758944
$arr[20,8] = 0.0787408966502241
$arr[20,9] = 0.4494545832162089
$arr[20,10] = 0.4786792303147678
$arr[20,11] = 0.281905505438516
$arr[20,12] = 1.99822666534058
$arr[20,13] = 3.777625341164565
$arr[21,0] = 0.7366151853966585
$arr[21,1] = 0.1437833108486188
$arr[21,2] = 0
$arr[21,3] = 0.2181108951600805
$arr[21,4] = 2.111450700512052
$arr[21,5] = 0.002471090746758872
$arr[21,6] = 0
$arr[21,7] = 1.001091277621988
$arr[21,8] = 0.07787989767514603
$arr[21,9] = 0.4352584355517308
$arr[21,10] = 0.4757086009489058
$arr[21,11] = 0.2779172836661559
$arr[21,12] = 2.005292204770345
$arr[21,13] = 3.783396002077012
$arr[22,0] = 0.6737658931359078
$arr[22,1] = 0.1425242454564639
$arr[22,2] = 0
$arr[22,3] = 0.2170846594251472
$arr[22,4] = 2.10491761703318
$arr[22,5] = 0.002474250227531713
$arr[22,6] = 0
$arr[22,7] = 1.009965054162901
$arr[22,8] = 0.07459601344045552
$arr[22,9] = 0.3815151725561918
$arr[22,10] = 0.4647448437439863
$arr[22,11] = 0.2629587066804291
$arr[22,12] = 2.033075746383667
$arr[22,13] = 3.808175995327048
$arr[23,0] = 0.6065651612016438
$arr[23,1] = 0.1411506279270327
$arr[23,2] = 0
$arr[23,3] = 0.2162303768465215
$arr[23,4] = 2.101295476114259
$arr[23,5] = 0.002477918982910064
$arr[23,6] = 0
$arr[23,7] = 1.021050016377664
$arr[23,8] = 0.07100543503205614
$arr[23,9] = 0.3236105103763975
$arr[23,10] = 0.4535385347643768
$arr[23,11] = 0.2471425519417494
$arr[23,12] = 2.065232390331913
$arr[23,13] = 3.84116056088169

$range = $ws.Range("B2:O25")
$range.Value = $arr